$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New submission row appended by the SmartScore Streamlit app.
$row = 25

# The SmartScore numeric-looking values in this particular export row were
# written out as plain text (not numbers), so format those cells as Text
# first to stop Excel from auto-coercing them to numbers on assignment.
$scoreCols = @("I", "L", "O", "R", "U", "X", "AA", "AD", "AG")
foreach ($col in $scoreCols) {
    $ws.Range($col + $row).NumberFormat = "@"
}

$ws.Range("A25").Value = "remas ali almadani_20251202_134128"
$ws.Range("B25").Value = ""
$ws.Range("C25").Value = "remas ali almadani"
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = "Female"
$ws.Range("F25").Value = "2025-12-02 13:41:28"
$ws.Range("G25").Value = @"
{
  "portion": 0.2,
  "diet": 0.2857142857142857,
  "salt": 0.6,
  "fat": 0.6,
  "natural": 0.8,
  "convenience": 0.4,
  "price": 1.0
}
"@
$ws.Range("H25").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("I25").Value = "0.578"
$ws.Range("J25").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
$ws.Range("K25").Value = "Maruchan Ramen Sabor Pollo"
$ws.Range("L25").Value = "0.566"
$ws.Range("M25").Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"
$ws.Range("N25").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("O25").Value = "0.455"
$ws.Range("P25").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"
$ws.Range("Q25").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("R25").Value = "0.712"
$ws.Range("S25").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
$ws.Range("T25").Value = "Annie’s Shells & White Cheddar"
$ws.Range("U25").Value = "0.625"
$ws.Range("V25").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"
$ws.Range("W25").Value = "Velveeta Original Shells & Cheese (microwave cups)"
$ws.Range("X25").Value = "0.567"
$ws.Range("Y25").Value = "Muy cremoso, porción individual, rápido, salado, ideal para niños"
$ws.Range("Z25").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AA25").Value = "0.657"
$ws.Range("AB25").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"
$ws.Range("AC25").Value = "Jack Link’s Beef Jerky Original"
$ws.Range("AD25").Value = "0.656"
$ws.Range("AE25").Value = "Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña"
$ws.Range("AF25").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AG25").Value = "0.644"
$ws.Range("AH25").Value = "Portátil, saludable, fácil, buena textura, sabor suave"

# Setting the multi-line JSON (the "Pesos" column) makes the host auto-expand
# the row height; AutoFit it back down so the row matches the sheet default
# (no explicit/custom row height), same as every other data row.
$ws.Rows($row).EntireRow.AutoFit()
